$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Name = "manipulations"
